$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in rows 2-4 (B1 stays the same)
$ws.Range("B2").Value = 373

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 163

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 94

# Row 5 (A5=3, B5=80) is removed entirely - clear its contents
$ws.Range("A5:B5").Clear()
